{"js": "// Word Office.js (Word JavaScript API) edit script.\n// Replaces the date heading and the 26 multiplication problems with\n// their updated values, matching the target diff exactly.\nconst replacements = [\n  [\"2023-09-10 Sunday\", \"2023-09-11 Monday\"],\n  [\"34\\u00d739=\", \"64\\u00d789=\"],\n  [\"47\\u00d760=\", \"13\\u00d753=\"],\n  [\"96\\u00d799=\", \"31\\u00d734=\"],\n  [\"82\\u00d791=\", \"59\\u00d736=\"],\n  [\"97\\u00d787=\", \"25\\u00d731=\"],\n  [\"51\\u00d767=\", \"99\\u00d716=\"],\n  [\"63\\u00d762=\", \"22\\u00d777=\"],\n  [\"46\\u00d744=\", \"44\\u00d749=\"],\n  [\"61\\u00d776=\", \"92\\u00d766=\"],\n  [\"23\\u00d744=\", \"67\\u00d761=\"],\n  [\"49\\u00d767=\", \"34\\u00d741=\"],\n  [\"71\\u00d726=\", \"17\\u00d723=\"],\n  [\"51\\u00d740=\", \"60\\u00d729=\"],\n  [\"59\\u00d790=\", \"21\\u00d777=\"],\n  [\"38\\u00d751=\", \"70\\u00d783=\"],\n  [\"72\\u00d772=\", \"62\\u00d763=\"],\n  [\"62\\u00d736=\", \"82\\u00d755=\"],\n  [\"21\\u00d712=\", \"33\\u00d763=\"],\n  [\"17\\u00d733=\", \"36\\u00d728=\"],\n  [\"80\\u00d712=\", \"27\\u00d739=\"],\n  [\"54\\u00d796=\", \"73\\u00d715=\"],\n  [\"80\\u00d745=\", \"75\\u00d733=\"],\n  [\"90\\u00d722=\", \"50\\u00d733=\"],\n  [\"38\\u00d762=\", \"52\\u00d790=\"],\n  [\"23\\u00d724=\", \"15\\u00d755=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Replaces the date heading and the 26 multiplication problems with\n# their updated values, matching the target diff exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-10 Sunday\", \"2023-09-11 Monday\"),\n    @(\"34\u00d739=\", \"64\u00d789=\"),\n    @(\"47\u00d760=\", \"13\u00d753=\"),\n    @(\"96\u00d799=\", \"31\u00d734=\"),\n    @(\"82\u00d791=\", \"59\u00d736=\"),\n    @(\"97\u00d787=\", \"25\u00d731=\"),\n    @(\"51\u00d767=\", \"99\u00d716=\"),\n    @(\"63\u00d762=\", \"22\u00d777=\"),\n    @(\"46\u00d744=\", \"44\u00d749=\"),\n    @(\"61\u00d776=\", \"92\u00d766=\"),\n    @(\"23\u00d744=\", \"67\u00d761=\"),\n    @(\"49\u00d767=\", \"34\u00d741=\"),\n    @(\"71\u00d726=\", \"17\u00d723=\"),\n    @(\"51\u00d740=\", \"60\u00d729=\"),\n    @(\"59\u00d790=\", \"21\u00d777=\"),\n    @(\"38\u00d751=\", \"70\u00d783=\"),\n    @(\"72\u00d772=\", \"62\u00d763=\"),\n    @(\"62\u00d736=\", \"82\u00d755=\"),\n    @(\"21\u00d712=\", \"33\u00d763=\"),\n    @(\"17\u00d733=\", \"36\u00d728=\"),\n    @(\"80\u00d712=\", \"27\u00d739=\"),\n    @(\"54\u00d796=\", \"73\u00d715=\"),\n    @(\"80\u00d745=\", \"75\u00d733=\"),\n    @(\"90\u00d722=\", \"50\u00d733=\"),\n    @(\"38\u00d762=\", \"52\u00d790=\"),\n    @(\"23\u00d724=\", \"15\u00d755=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdFindContinue = 1 (Wrap), wdReplaceAll = 2 (Replace)\n    $found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
